# Update Facebook and Twitter data: revise existing quarterly counts
# and append three new quarters (2021-12-31, 2022-03-31, 2022-06-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing counts in column B (rows 2-12) ---
$ws.Cells.Item(2, 2).Value = 1796
$ws.Cells.Item(3, 2).Value = 1816
$ws.Cells.Item(4, 2).Value = 1766
$ws.Cells.Item(5, 2).Value = 1734
$ws.Cells.Item(6, 2).Value = 1740
$ws.Cells.Item(7, 2).Value = 1670
$ws.Cells.Item(8, 2).Value = 1653
$ws.Cells.Item(9, 2).Value = 1670
$ws.Cells.Item(10, 2).Value = 1607
$ws.Cells.Item(11, 2).Value = 1669
$ws.Cells.Item(12, 2).Value = 1645

# --- Append new rows 13-15 ---
$newRows = @(
    @{ Row = 13; Date = 44561; Count = 1598 },
    @{ Row = 14; Date = 44651; Count = 1617 },
    @{ Row = 15; Date = 44742; Count = 1588 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r.Row, 2).Value = $r.Count
    $ws.Cells.Item($r.Row, 3).Value = "Q"
}
